# Auto-generated edit script: updates cryptos list values (prices / 1h
# volume deltas) + fixes a swapped NEARProtocol/PancakeSwap row pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.154.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '''3.070.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''587.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '''151.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.65%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.545'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.84%  '
$ws.Range("D9").Value = '''3.062.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("D11").Value = '''5.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '''0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '''0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").Value = '''37.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '''3.579.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").Value = '''7.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '''63.258.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '''3.069.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").Value = '''473.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = '''14.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.59%  '
$ws.Range("D22").Value = '''0.715'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").Value = '''7.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").Value = '''2.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").Value = '''13.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '''81.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D28").Value = '''9.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.07%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '''7.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '''2.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("D33").Value = '''0.114'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.16%  '
$ws.Range("D34").Value = '''27.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '''0.0₃0842'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '''6.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").Value = '''3.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("D39").Value = '''2.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.19%  '
$ws.Range("D40").Value = '''9.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("D41").Value = '''50.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.07%  '
$ws.Range("D42").Value = '''443.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.91%  '
$ws.Range("D43").Value = '''0.285'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.84%  '
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("D45").Value = '''39.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.38%  '
$ws.Range("D46").Value = '''2.805.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").Value = '''131.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D50").Value = '''25.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.87%  '
$ws.Range("D51").Value = '''2.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
